# Update cryptocurrency Price (column D) and Volume(1h) (column E) values
# on the active worksheet to reflect the latest scrape, per the
# "Updated symbol list" GitHub Actions commit.
#
# The Price/Volume columns store text values (e.g. "303.35", "12.47%"),
# so each cell's number format is forced to Text ("@") before assigning
# the new string value. This prevents Excel from re-interpreting the
# text as a number (which would also destroy trailing zeros/precision).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "303.35"),
    @("D3", "34.90"),
    @("E3", "12.47%"),
    @("D4", "5.166"),
    @("D5", "0.07828"),
    @("E5", "6.38%"),
    @("D6", "2.337"),
    @("E6", "2.42%"),
    @("D7", "8.054"),
    @("E7", "4.26%"),
    @("D8", "3.972"),
    @("E8", "6.37%"),
    @("D9", "0.9249"),
    @("E9", "0.83%"),
    @("D10", "0.1006"),
    @("E10", "8.55%"),
    @("D11", "0.1831"),
    @("E11", "7.45%"),
    @("D12", "0.08520"),
    @("E12", "2.71%"),
    @("D13", "0.03433"),
    @("E13", "10.54%"),
    @("D14", "0.09901"),
    @("E14", "-0.77%"),
    @("D15", "0.001470"),
    @("E15", "-1.64%"),
    @("D16", "0.005764"),
    @("E16", "-0.16%"),
    @("D17", "3.472"),
    @("E17", "0.07%"),
    @("E18", "3.25%"),
    @("D19", "0.3427"),
    @("E19", "3.13%"),
    @("D20", "0.1327"),
    @("E20", "2.98%"),
    @("D21", "4.556"),
    @("E21", "9.62%"),
    @("D22", "0.2223"),
    @("E22", "4.70%"),
    @("D23", "0.04652"),
    @("E23", "3.03%"),
    @("D24", "0.001219"),
    @("E24", "0.70%"),
    @("D25", "0.004453"),
    @("E25", "6.40%"),
    @("D26", "0.0001300"),
    @("E26", "0.09%"),
    @("D27", "0.0003399"),
    @("E27", "0.21%"),
    @("D39", "0.01751"),
    @("E39", "11.14%"),
    @("D40", "0.04751"),
    @("E40", "5.82%"),
    @("D41", "0.007763"),
    @("E41", "5.59%"),
    @("E42", "5.73%"),
    @("D43", "0.008854"),
    @("E43", "-9.96%"),
    @("D44", "0.002211"),
    @("E44", "-0.81%"),
    @("D45", "0.009974"),
    @("E45", "7.92%"),
    @("D46", "0.00006075"),
    @("E46", "-0.35%"),
    @("D47", "0.00000000750"),
    @("E47", "0.10%"),
    @("D48", "5.837"),
    @("E48", "123.98%"),
    @("D49", "0.002691"),
    @("E49", "28.22%"),
    @("D50", "0.00002100"),
    @("E50", "0.10%"),
    @("D51", "0.0002000"),
    @("E51", "0.10%")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newVal
}
